$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.878.92"
$ws.Range("E2").Value = "  -0.72%  "
$ws.Range("D3").Value = "3.453.60"
$ws.Range("E3").Value = "  -1.17%  "
$ws.Range("D5").Formula = "'574.74"
$ws.Range("E5").Value = "  -1.18%  "
$ws.Range("D6").Formula = "'159.68"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").Value = "3.452.14"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").Formula = "'0.574"
$ws.Range("E9").Value = "  -6.25%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  -3.19%  "
$ws.Range("D12").Formula = "'0.440"
$ws.Range("E12").Value = "  -1.82%  "
$ws.Range("D13").Value = "4.048.69"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("D15").Formula = "'27.64"
$ws.Range("E15").Value = "  -4.10%  "
$ws.Range("D16").Formula = "'0.0000174"
$ws.Range("E16").Value = "  -10.25%  "
$ws.Range("D17").Value = "64.924.99"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("D18").Value = "3.430.24"
$ws.Range("E18").Value = "  -1.77%  "
$ws.Range("D19").Formula = "'6.21"
$ws.Range("E19").Value = "  -4.03%  "
$ws.Range("D20").Formula = "'13.74"
$ws.Range("E20").Value = "  -4.76%  "
$ws.Range("D21").Formula = "'377.76"
$ws.Range("E21").Value = "  -2.11%  "
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").Formula = "'72.31"
$ws.Range("E24").Value = "  -0.60%  "
$ws.Range("E25").Value = "  -3.41%  "
$ws.Range("E26").Value = "  -0.63%  "
$ws.Range("D27").Formula = "'9.91"
$ws.Range("E27").Value = "  -2.40%  "
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Formula = "'1.00"
$ws.Range("E29").Value = "  +0.08%  "
$ws.Range("E30").Value = "  -5.31%  "
$ws.Range("D31").Formula = "'6.06"
$ws.Range("E31").Value = "  -2.04%  "
$ws.Range("E32").Value = "  -2.41%  "
$ws.Range("D33").Formula = "'23.17"
$ws.Range("E33").Value = "  -2.43%  "
$ws.Range("D34").Formula = "'6.99"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").Formula = "'1.56"
$ws.Range("E35").Value = "  -4.19%  "
$ws.Range("D36").Formula = "'161.12"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Formula = "'1.87"
$ws.Range("E37").Value = "  -3.32%  "
$ws.Range("D38").Value = "2.900.43"
$ws.Range("E38").Value = "  -4.43%  "
$ws.Range("D39").Formula = "'0.0749"
$ws.Range("E39").Value = "  -4.20%  "
$ws.Range("D40").Formula = "'26.18"
$ws.Range("E40").Value = "  -3.31%  "
$ws.Range("D41").Formula = "'4.52"
$ws.Range("E41").Value = "  -2.00%  "
$ws.Range("D42").Formula = "'43.02"
$ws.Range("E42").Value = "  -0.04%  "
$ws.Range("D43").Formula = "'6.51"
$ws.Range("E43").Value = "  -5.56%  "
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Formula = "'26.20"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("D46").Formula = "'0.0311"
$ws.Range("E46").Value = "  -3.42%  "
$ws.Range("D47").Formula = "'2.37"
$ws.Range("E47").Value = "  +7.57%  "
$ws.Range("D48").Formula = "'321.24"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("E49").Value = "  -2.92%  "
$ws.Range("D50").Formula = "'6.47"
$ws.Range("E50").Value = "  -4.35%  "
$ws.Range("E51").Value = "  -4.09%  "
